# Fix "Время начала заседания: [time]" placeholder -> "Время начала заседания: _____"
# The original paragraph spelled out the placeholder as an underlined "[time]"
# (split across several runs, including spell-check proofErr markers around the
# word "time"). Replace it with a plain underscore blank, matching the
# surrounding (non-underlined) run formatting.
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Время начала заседания: [time]",  # FindText
    $true,                             # MatchCase
    $false,                            # MatchWholeWord
    $false,                            # MatchWildcards (brackets are literal)
    $false,                            # MatchSoundsLike
    $false,                            # MatchAllWordForms
    $true,                             # Forward
    1,                                 # Wrap -> wdFindContinue
    $false,                            # Format
    "Время начала заседания: _____",   # ReplaceWith
    2                                  # Replace -> wdReplaceAll
)

Write-Output "Replaced placeholder: $found"
